# Apply "Updated symbol list" edits to the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds numeric-looking values that must remain
# stored as text (matching the original inlineStr cells, e.g. "248.30"
# keeps its trailing zero). Force text formatting on the column first so
# Excel does not silently convert the assigned strings into numbers.
$ws.Range("D2:D49").NumberFormat = "@"

# --- Price (column D) updates ---
$ws.Cells.Item(2, 4).Value = "248.30"
$ws.Cells.Item(3, 4).Value = "21.79"
$ws.Cells.Item(4, 4).Value = "5.498"
$ws.Cells.Item(5, 4).Value = "0.05646"
$ws.Cells.Item(6, 4).Value = "3.378"
$ws.Cells.Item(7, 4).Value = "6.441"
$ws.Cells.Item(9, 4).Value = "1.040"
$ws.Cells.Item(11, 4).Value = "0.07244"
$ws.Cells.Item(12, 4).Value = "0.03111"
$ws.Cells.Item(13, 4).Value = "0.02940"
$ws.Cells.Item(14, 4).Value = "0.09290"
$ws.Cells.Item(15, 4).Value = "0.001647"
$ws.Cells.Item(16, 4).Value = "3.211"
$ws.Cells.Item(17, 4).Value = "0.04718"
$ws.Cells.Item(18, 4).Value = "0.0005826"
$ws.Cells.Item(19, 4).Value = "0.006440"
$ws.Cells.Item(20, 4).Value = "0.005030"
$ws.Cells.Item(21, 4).Value = "0.001049"
$ws.Cells.Item(24, 4).Value = "4.162"
$ws.Cells.Item(40, 4).Value = "0.04094"
$ws.Cells.Item(42, 4).Value = "0.003502"
$ws.Cells.Item(44, 4).Value = "0.008457"
$ws.Cells.Item(45, 4).Value = "0.00005813"
$ws.Cells.Item(46, 4).Value = "0.00000000750"
$ws.Cells.Item(47, 4).Value = "0.7859"
$ws.Cells.Item(49, 4).Value = "0.00002101"

# --- Row 18: "17OneONE" -> "17OneONEWorstin24h" ---
$ws.Cells.Item(18, 5).Value = "17OneONEWorstin24h"

# --- Row 41: BKEXToken -> KickToken ---
$ws.Cells.Item(41, 2).Value = "KickToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Cells.Item(41, 4).Value = "0.006935"
$ws.Cells.Item(41, 5).Value = "40KickTokenKICK"

# --- Row 43: KickToken -> BKEXToken ---
$ws.Cells.Item(43, 2).Value = "BKEXToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Cells.Item(43, 4).Value = "0.1041"
$ws.Cells.Item(43, 5).Value = "42BKEXTokenBKK"
